$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 23, shifting existing rows 23-28 down to 25-30
$ws.Rows.Item(23).Resize(2).Insert()

# Row 23: norway / lau2 / position_geolabels / 2020 / default / data.table
$ws.Range("A23").Value = "norway_lau2_position_geolabels_b2020_default_dt"
$ws.Range("B23").Value = "norway"
$ws.Range("C23").Value = "lau2"
$ws.Range("D23").Value = "position_geolabels"
$ws.Range("E23").Value = 2020
$ws.Range("F23").Value = "default"
$ws.Range("G23").Value = "data.table"

# Row 24: norway / lau2 / position_geolabels / 2020 / insert_oslo / data.table
$ws.Range("A24").Value = "norway_lau2_position_geolabels_b2020_insert_oslo_dt"
$ws.Range("B24").Value = "norway"
$ws.Range("C24").Value = "lau2"
$ws.Range("D24").Value = "position_geolabels"
$ws.Range("E24").Value = 2020
$ws.Range("F24").Value = "insert_oslo"
$ws.Range("G24").Value = "data.table"

# Update selection to match the final state
$ws.Range("A25").Select()
